$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 107
$ws.Range("H107").Value = 657
$ws.Range("I107").Value = 599.8333
$ws.Range("K107").Value = 599.8333
$ws.Range("M107").Value = 1320.1667
# Row 111
$ws.Range("H111").Value = 2921.5715
$ws.Range("I111").Value = 2683.8
$ws.Range("J111").Value = 3516
$ws.Range("K111").Value = 8051.400000000001
$ws.Range("L111").Value = 10548
$ws.Range("M111").Value = -4984.400000000001
$ws.Range("N111").Value = -16682
# Row 112
$ws.Range("H112").Value = 1972.6666
$ws.Range("J112").Value = 1972.6666
$ws.Range("L112").Value = 5917.9998
$ws.Range("N112").Value = -8133.9998
# Row 113
$ws.Range("H113").Value = 2110
$ws.Range("I113").Value = 2013.3334
$ws.Range("J113").Value = 2980
$ws.Range("K113").Value = 2013.3334
$ws.Range("L113").Value = 2980
$ws.Range("M113").Value = 1240.6666
$ws.Range("N113").Value = -9488
# Row 115
$ws.Range("H115").Value = 0
$ws.Range("I115").Value = 0
$ws.Range("K115").Value = 0
$ws.Range("M115").ClearContents()
# Row 116
$ws.Range("H116").Value = 3542.5386
$ws.Range("I116").Value = 4612.1177
$ws.Range("J116").Value = 1522.2222
$ws.Range("K116").Value = 4612.1177
$ws.Range("L116").Value = 1522.2222
$ws.Range("M116").Value = -1170.1177
$ws.Range("N116").Value = -8406.2222
# Row 118
$ws.Range("H118").Value = 1717.6471
$ws.Range("I118").Value = 900
$ws.Range("J118").Value = 1768.75
$ws.Range("K118").Value = 2700
$ws.Range("L118").Value = 5306.25
$ws.Range("M118").Value = -1043
$ws.Range("N118").Value = -8620.25
# Row 130
$ws.Range("H130").Value = 44422.855
$ws.Range("J130").Value = 44422.855
$ws.Range("L130").Value = 44422.855
$ws.Range("N130").Value = -54462.855
# Row 137
$ws.Range("H137").Value = 5310.759
$ws.Range("I137").Value = 1434.3334
$ws.Range("J137").Value = 6322
$ws.Range("K137").Value = 4303.0002
$ws.Range("L137").Value = 18966
$ws.Range("M137").Value = -1753.0002
$ws.Range("N137").Value = -24066

$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 2052.5
$ws.Range("I2").Value = 2211.2222
$ws.Range("J2").Value = 1848.4286
$ws.Range("K2").Value = 2211.2222
$ws.Range("L2").Value = 1848.4286
$ws.Range("M2").Value = -2098.2222
$ws.Range("N2").Value = -2074.4286
# Row 74
$ws.Range("H74").Value = 4041.5454
$ws.Range("I74").Value = 882.4583
$ws.Range("J74").Value = 12465.777
$ws.Range("K74").Value = 882.4583
$ws.Range("L74").Value = 12465.777
$ws.Range("M74").Value = -8.458300000000008
$ws.Range("N74").Value = -14213.777
# Row 77
$ws.Range("H77").Value = 4041.5454
$ws.Range("I77").Value = 882.4583
$ws.Range("J77").Value = 12465.777
$ws.Range("K77").Value = 4412.2915
$ws.Range("L77").Value = 62328.885
$ws.Range("M77").Value = -44.29150000000027
$ws.Range("N77").Value = -71064.88500000001
# Row 82
$ws.Range("H82").Value = 40200
$ws.Range("J82").Value = 40200
$ws.Range("L82").Value = 40200
$ws.Range("N82").Value = -40922
# Row 85
$ws.Range("H85").Value = 40200
$ws.Range("J85").Value = 40200
$ws.Range("L85").Value = 40200
$ws.Range("N85").Value = -42696
# Row 110
$ws.Range("H110").Value = 1732.3704
$ws.Range("I110").Value = 1738
$ws.Range("J110").Value = 1700
$ws.Range("K110").Value = 1738
$ws.Range("L110").Value = 1700
$ws.Range("M110").Value = 307
$ws.Range("N110").Value = -5790
# Row 116
$ws.Range("H116").Value = 2052.5
$ws.Range("I116").Value = 2211.2222
$ws.Range("J116").Value = 1848.4286
$ws.Range("K116").Value = 2211.2222
$ws.Range("L116").Value = 1848.4286
$ws.Range("M116").Value = 82.77779999999984
$ws.Range("N116").Value = -6436.4286

$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 2052.5
$ws.Range("I3").Value = 2211.2222
$ws.Range("J3").Value = 1848.4286
$ws.Range("K3").Value = 2211.2222
$ws.Range("L3").Value = 1848.4286
$ws.Range("M3").Value = -2097.2222
$ws.Range("N3").Value = -2076.4286
# Row 105
$ws.Range("H105").Value = 3876.7896
$ws.Range("I105").Value = 7025
$ws.Range("K105").Value = 7025
$ws.Range("M105").Value = -5278
# Row 107
$ws.Range("H107").Value = 1769.6666
$ws.Range("I107").Value = 1454.75
$ws.Range("J107").Value = 2189.5557
$ws.Range("K107").Value = 1454.75
$ws.Range("L107").Value = 2189.5557
$ws.Range("M107").Value = 465.25
$ws.Range("N107").Value = -6029.5557

$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 9166.666999999999
$ws.Range("I16").Value = 1500
$ws.Range("K16").Value = 1500
$ws.Range("M16").Value = -1213
# Row 31
$ws.Range("H31").Value = 6045.3145
$ws.Range("I31").Value = 1789.8572
$ws.Range("J31").Value = 8882.286
$ws.Range("K31").Value = 1789.8572
$ws.Range("L31").Value = 8882.286
$ws.Range("M31").Value = -1494.8572
$ws.Range("N31").Value = -9472.286
# Row 34
$ws.Range("H34").Value = 6045.3145
$ws.Range("I34").Value = 1789.8572
$ws.Range("J34").Value = 8882.286
$ws.Range("K34").Value = 1789.8572
$ws.Range("L34").Value = 8882.286
$ws.Range("M34").Value = -1587.8572
$ws.Range("N34").Value = -9286.286
# Row 105
$ws.Range("H105").Value = 1048.8572
$ws.Range("I105").Value = 985.625
$ws.Range("J105").Value = 1133.1666
$ws.Range("K105").Value = 985.625
$ws.Range("L105").Value = 1133.1666
$ws.Range("M105").Value = 761.375
$ws.Range("N105").Value = -4627.1666
# Row 113
$ws.Range("H113").Value = 9166.666999999999
$ws.Range("I113").Value = 1500
$ws.Range("K113").Value = 1500
$ws.Range("M113").Value = 670

$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 495.77966
$ws.Range("I5").Value = 429.1698
$ws.Range("J5").Value = 1084.1666
$ws.Range("K5").Value = 1287.5094
$ws.Range("L5").Value = 3252.4998
$ws.Range("M5").Value = -1175.5094
$ws.Range("N5").Value = -3476.4998
# Row 108
$ws.Range("H108").Value = 2165
$ws.Range("I108").Value = 400
$ws.Range("J108").Value = 3930
$ws.Range("K108").Value = 1200
$ws.Range("L108").Value = 11790
$ws.Range("M108").Value = 1680
$ws.Range("N108").Value = -17550
# Row 115
$ws.Range("H115").Value = 1166.6666
$ws.Range("I115").Value = 1166.6666
$ws.Range("K115").Value = 3499.9998
$ws.Range("M115").Value = -2324.9998
# Row 122
$ws.Range("H122").Value = 884.6316
$ws.Range("I122").Value = 432.8
$ws.Range("J122").Value = 1046
$ws.Range("K122").Value = 3895.2
$ws.Range("L122").Value = 9414
$ws.Range("M122").Value = -1445.2
$ws.Range("N122").Value = -14314
# Row 135
$ws.Range("H135").Value = 495.77966
$ws.Range("I135").Value = 429.1698
$ws.Range("J135").Value = 1084.1666
$ws.Range("K135").Value = 3862.5282
$ws.Range("L135").Value = 9757.499400000001
$ws.Range("M135").Value = -1327.5282
$ws.Range("N135").Value = -14827.4994

$ws = $wb.Worksheets.Item("GSM")
# Row 107
$ws.Range("H107").Value = 250
$ws.Range("I107").Value = 233.33333
$ws.Range("J107").Value = 300
$ws.Range("K107").Value = 233.33333
$ws.Range("L107").Value = 300
$ws.Range("M107").Value = 1686.66667
$ws.Range("N107").Value = -4140
# Row 113
$ws.Range("H113").Value = 4275.421
$ws.Range("I113").Value = 3661.0833
$ws.Range("J113").Value = 5328.5713
$ws.Range("K113").Value = 3661.0833
$ws.Range("L113").Value = 5328.5713
$ws.Range("M113").Value = -1491.0833
$ws.Range("N113").Value = -9668.5713

$ws = $wb.Worksheets.Item("LTW")
# Row 40
$ws.Range("H40").Value = 2894.0454
$ws.Range("I40").Value = 2695.394
$ws.Range("K40").Value = 2695.394
$ws.Range("M40").Value = -2559.394
# Row 46
$ws.Range("H46").Value = 1535.9166
$ws.Range("J46").Value = 875
$ws.Range("L46").Value = 875
$ws.Range("N46").Value = -1251
# Row 61
$ws.Range("H61").Value = 2181.8462
$ws.Range("J61").Value = 2560
$ws.Range("L61").Value = 2560
$ws.Range("N61").Value = -2964
# Row 93
$ws.Range("H93").Value = 2386.5908
$ws.Range("I93").Value = 2377.3076
$ws.Range("J93").Value = 2400
$ws.Range("K93").Value = 2377.3076
$ws.Range("L93").Value = 2400
$ws.Range("M93").Value = -1129.3076
$ws.Range("N93").Value = -4896
# Row 100
$ws.Range("H100").Value = 1245.7778
$ws.Range("I100").Value = 1144.5714
$ws.Range("J100").Value = 1600
$ws.Range("K100").Value = 1144.5714
$ws.Range("L100").Value = 1600
$ws.Range("M100").Value = -603.5714
$ws.Range("N100").Value = -2682
# Row 113
$ws.Range("H113").Value = 2181.8462
$ws.Range("J113").Value = 2560
$ws.Range("L113").Value = 2560
$ws.Range("N113").Value = -6900
# Row 122
$ws.Range("H122").Value = 62239.234
$ws.Range("I122").Value = 86379.75
$ws.Range("K122").Value = 259139.25
$ws.Range("M122").Value = -256689.25

$ws = $wb.Worksheets.Item("WVR")
# Row 94
$ws.Range("H94").Value = 43999.5
$ws.Range("J94").Value = 43999.5
$ws.Range("L94").Value = 43999.5
$ws.Range("N94").Value = -45801.5
# Row 136
$ws.Range("H136").Value = 4186.6665
$ws.Range("I136").Value = 5980.263
$ws.Range("J136").Value = 2182.0588
$ws.Range("K136").Value = 17940.789
$ws.Range("L136").Value = 6546.176399999999
$ws.Range("M136").Value = -15390.789
$ws.Range("N136").Value = -11646.1764
